$d = $word.ActiveDocument

$replacements = @(
    @("2025-09-03 Wednesday", "2025-09-04 Thursday"),
    @("95×94=8930", "70×13=910"),
    @("83×43=3569", "66×41=2706"),
    @("11×50=550", "60×90=5400"),
    @("83×22=1826", "79×16=1264"),
    @("97×85=8245", "13×88=1144"),
    @("52×73=3796", "71×62=4402"),
    @("75×92=6900", "30×32=960"),
    @("12×97=1164", "30×56=1680"),
    @("11×90=990", "94×18=1692"),
    @("29×68=1972", "32×53=1696"),
    @("97×94=9118", "58×82=4756"),
    @("52×23=1196", "68×55=3740"),
    @("29×83=2407", "19×83=1577"),
    @("21×14=294", "68×31=2108"),
    @("65×98=6370", "32×44=1408"),
    @("13×11=143", "30×91=2730"),
    @("39×61=2379", "17×89=1513"),
    @("50×32=1600", "25×91=2275"),
    @("53×31=1643", "61×67=4087"),
    @("68×14=952", "88×70=6160"),
    @("22×85=1870", "93×40=3720"),
    @("97×76=7372", "21×15=315"),
    @("87×39=3393", "73×74=5402"),
    @("82×46=3772", "81×98=7938"),
    @("66×81=5346", "42×46=1932")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
